$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: amount / invoice no / customer / description
$ws.Range("C2").Value = 30000
$ws.Range("D2").Value = 1109
$ws.Range("E2").Value = "WALMART"
$ws.Range("F2").Value = "Tenth Order"

# Row 3: amount / invoice no / customer / description
$ws.Range("C3").Value = 35000
$ws.Range("D3").Value = 1106
$ws.Range("E4").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "HOMEDEP"
$ws.Range("F3").Value = "Tenth Order"

# Row 4 becomes an empty (but still formatted) row
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

# Row 5 is removed entirely, shifting nothing below it up
$ws.Range("A5:F5").Delete()
